$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: swap Agence/Banque values (D2/E2), and clear Taxe/MT Net (J2/K2)
$ws.Range("D2").Value = "BMCE MAARIF"
$ws.Range("E2").Value = "BMCE"
$ws.Range("J2").Value = ""
$ws.Range("K2").Value = ""

# Row 4: swap Agence/Banque values (D4/E4)
$ws.Range("D4").Value = "BMCE test"
$ws.Range("E4").Value = "BMCE"

# Row 5: swap Agence/Banque values (D5/E5)
$ws.Range("D5").Value = ""
$ws.Range("E5").Value = "12121"

# Row 6: new entry - Hassan Hssouni
$ws.Range("A6").Value = "Hassan Hssouni"
$ws.Range("B6").Value = "BJ123456"
$ws.Range("C6").Value = "116497823245768736541324"
$ws.Range("D6").Value = "BMCE TARIQ"
$ws.Range("E6").Value = "BMCE"
$ws.Range("F6").Value = "Supervision"
$ws.Range("G6").Value = "123/Test SUP"
$ws.Range("H6").Value = "trimestrielle"
$ws.Range("I6").Value = 100000
$ws.Range("J6").Value = 22500
$ws.Range("K6").Value = 92500

# Row 7: new entry - Amine Kamal
$ws.Range("A7").Value = "Amine Kamal"
$ws.Range("B7").Value = "cd1200"
$ws.Range("C7").Value = "022232265645652220000000"
$ws.Range("D7").Value = "almoqawama"
$ws.Range("E7").Value = "bmce"
$ws.Range("F7").Value = "Direction régionale"
$ws.Range("G7").Value = "002/DR002"
$ws.Range("H7").Value = "mensuelle"
$ws.Range("I7").Value = 40000
$ws.Range("J7").Value = 3000
$ws.Range("K7").Value = 37000
